# Updates cryptos list with latest prices / 1h volume changes
# (GitHub Actions scheduled refresh)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, never letting Excel
# auto-convert numeric-looking strings (e.g. "1.010") into numbers.
# Number-format the cell as Text, assign, then restore the cell to the
# workbook's default (General / Normal) so no stray style is left behind.
function Set-TextValue($addr, [string]$val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "28.291.36"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "1.887.55"
$ws.Range("E3").Value = "  +1.44%  "
Set-TextValue "D4" "1.010"
$ws.Range("E4").Value = "  +0.69%  "
Set-TextValue "D5" "316.38"
$ws.Range("E5").Value = "  +1.32%  "
Set-TextValue "D6" "1.011"
$ws.Range("E6").Value = "  +0.81%  "
Set-TextValue "D7" "0.5147"
$ws.Range("E7").Value = "  +0.61%  "
Set-TextValue "D8" "0.3907"
$ws.Range("E8").Value = "  +1.36%  "
Set-TextValue "D9" "0.08398"
$ws.Range("E9").Value = "  +1.49%  "
Set-TextValue "D10" "1.126"
$ws.Range("E10").Value = "  +1.19%  "
Set-TextValue "D11" "41.72"
$ws.Range("E11").Value = "  +0.39%  "
Set-TextValue "D12" "6.252"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "1.879.94"
$ws.Range("E13").Value = "  +0.71%  "
Set-TextValue "D14" "20.62"
$ws.Range("E14").Value = "  +0.14%  "
Set-TextValue "D15" "7.280"
$ws.Range("E15").Value = "  +0.45%  "
Set-TextValue "D16" "1.010"
$ws.Range("E16").Value = "  +0.77%  "
Set-TextValue "D17" "0.00001108"
$ws.Range("E17").Value = "  +0.91%  "
Set-TextValue "D18" "91.17"
$ws.Range("E18").Value = "  +0.56%  "
Set-TextValue "D19" "0.06697"
$ws.Range("E19").Value = "  +0.55%  "
Set-TextValue "D20" "17.83"
$ws.Range("E20").Value = "  +0.73%  "
Set-TextValue "D21" "1.010"
$ws.Range("E21").Value = "  +0.79%  "
Set-TextValue "D22" "6.042"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "28.320.17"
$ws.Range("E23").Value = "  +1.08%  "
Set-TextValue "D24" "11.18"
Set-TextValue "D25" "2.293"
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("D26").Value = "2.101.72"
Set-TextValue "D27" "160.60"
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("E28").Value = "  -1.87%  "
Set-TextValue "D29" "20.73"
$ws.Range("E29").Value = "  +1.11%  "
Set-TextValue "D30" "125.68"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +1.08%  "
Set-TextValue "D33" "5.897"
$ws.Range("E33").Value = "  -0.76%  "
Set-TextValue "D34" "3.627"
$ws.Range("E34").Value = "  +0.92%  "
Set-TextValue "D35" "9.489"
$ws.Range("E35").Value = "  +1.14%  "
Set-TextValue "D36" "0.02449"
$ws.Range("E36").Value = "  +1.67%  "
Set-TextValue "D37" "0.06591"
$ws.Range("E37").Value = "  +1.48%  "
Set-TextValue "D38" "0.2213"
$ws.Range("E38").Value = "  +1.65%  "
Set-TextValue "D39" "1.201"
$ws.Range("E39").Value = "  +0.40%  "
Set-TextValue "D40" "0.6512"
$ws.Range("E40").Value = "  -1.70%  "
Set-TextValue "D41" "1.249"
$ws.Range("E41").Value = "  +2.28%  "
Set-TextValue "D42" "5.017"
$ws.Range("E42").Value = "  +0.10%  "
Set-TextValue "D43" "11.26"
$ws.Range("E43").Value = "  +0.71%  "
Set-TextValue "D44" "0.6111"
$ws.Range("E44").Value = "  -0.87%  "
Set-TextValue "D45" "13.14"
$ws.Range("E45").Value = "  +0.77%  "
Set-TextValue "D46" "3.702"
$ws.Range("E46").Value = "  +1.35%  "
Set-TextValue "D47" "1.286"
$ws.Range("E47").Value = "  +0.54%  "
Set-TextValue "D48" "2.020"
$ws.Range("E48").Value = "  +0.48%  "
Set-TextValue "D49" "1.239"
$ws.Range("E49").Value = "  +2.64%  "
Set-TextValue "D50" "121.34"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.06920"
$ws.Range("E51").Value = "  +1.25%  "
